# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Reorders/refreshes the worker period-of-arrears table (rows 16-29) so
# each worker's periods run chronologically (2102 -> 2108) and the two
# workers are interleaved row-by-row, with refreshed "Valor Mora" /
# "Salario Basico" figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row, TipoDoc, NumDoc, Nombre, Periodo, ValorMora, SalarioBasico
$rows = @(
    @(16, "CC", "73207105",   "DAVID HERNANDO SOTO NUÑEZ",      "2102", 120000, 3000000),
    @(17, "CC", "1143393622", "AURA PATRICIA PEÑARANDA SEGURA", "2102", 35112,  877803),
    @(18, "CC", "73207105",   "DAVID HERNANDO SOTO NUÑEZ",      "2103", 120000, 3000000),
    @(19, "CC", "1143393622", "AURA PATRICIA PEÑARANDA SEGURA", "2103", 35112,  877803),
    @(20, "CC", "73207105",   "DAVID HERNANDO SOTO NUÑEZ",      "2104", 120000, 3000000),
    @(21, "CC", "1143393622", "AURA PATRICIA PEÑARANDA SEGURA", "2104", 35112,  877803),
    @(22, "CC", "73207105",   "DAVID HERNANDO SOTO NUÑEZ",      "2105", 120000, 3000000),
    @(23, "CC", "1143393622", "AURA PATRICIA PEÑARANDA SEGURA", "2105", 35112,  877803),
    @(24, "CC", "73207105",   "DAVID HERNANDO SOTO NUÑEZ",      "2106", 120000, 3000000),
    @(25, "CC", "1143393622", "AURA PATRICIA PEÑARANDA SEGURA", "2106", 35112,  877803),
    @(26, "CC", "73207105",   "DAVID HERNANDO SOTO NUÑEZ",      "2107", 120000, 3000000),
    @(27, "CC", "1143393622", "AURA PATRICIA PEÑARANDA SEGURA", "2107", 35112,  877803),
    @(28, "CC", "73207105",   "DAVID HERNANDO SOTO NUÑEZ",      "2108", 116000, 3000000),
    @(29, "CC", "1143393622", "AURA PATRICIA PEÑARANDA SEGURA", "2108", 33942,  877803)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 2).Value = $r[1]
    $ws.Cells.Item($rowNum, 3).Value = $r[2]
    $ws.Cells.Item($rowNum, 4).Value = $r[3]
    $ws.Cells.Item($rowNum, 5).Value = $r[4]
    $ws.Cells.Item($rowNum, 6).Value = $r[5]
    $ws.Cells.Item($rowNum, 7).Value = $r[6]
}
